$ws = $excel.ActiveWorkbook.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

$sub3 = [System.Text.Encoding]::UTF8.GetString([byte[]](0xE2,0x82,0x83))
$sub6 = [System.Text.Encoding]::UTF8.GetString([byte[]](0xE2,0x82,0x86))
$d32val = "0.0" + $sub3 + "0892"
$d51val = "0.0" + $sub6 + "0233"

# Row 2
Set-TextValue $ws.Range("D2") "63.804.64"
Set-TextValue $ws.Range("E2") "  +0.14%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.620.57"
Set-TextValue $ws.Range("E3") "  +0.11%  "

# Row 4
Set-TextValue $ws.Range("E4") "  -0.05%  "

# Row 5
Set-TextValue $ws.Range("D5") "595.20"
Set-TextValue $ws.Range("E5") "  +0.06%  "

# Row 6
Set-TextValue $ws.Range("D6") "150.93"
Set-TextValue $ws.Range("E6") "  +1.05%  "

# Row 7
Set-TextValue $ws.Range("E7") "  -0.07%  "

# Row 8
Set-TextValue $ws.Range("E8") "  -0.15%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.114"
Set-TextValue $ws.Range("E9") "  +4.92%  "

# Row 10
Set-TextValue $ws.Range("D10") "5.81"
Set-TextValue $ws.Range("E10") "  +4.24%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.395"
Set-TextValue $ws.Range("E11") "  +3.68%  "

# Row 12
Set-TextValue $ws.Range("E12") "  +1.01%  "

# Row 13
Set-TextValue $ws.Range("D13") "27.94"
Set-TextValue $ws.Range("E13") "  +1.84%  "

# Row 14
Set-TextValue $ws.Range("D14") "3.088.72"
Set-TextValue $ws.Range("E14") "  -0.02%  "

# Row 15
Set-TextValue $ws.Range("D15") "63.592.23"
Set-TextValue $ws.Range("E15") "  -0.03%  "

# Row 16
Set-TextValue $ws.Range("D16") "0.0000164"
Set-TextValue $ws.Range("E16") "  +11.17%  "

# Row 17
Set-TextValue $ws.Range("D17") "2.594.10"
Set-TextValue $ws.Range("E17") "  -0.72%  "

# Row 18
Set-TextValue $ws.Range("D18") "12.25"
Set-TextValue $ws.Range("E18") "  +1.15%  "

# Row 19
Set-TextValue $ws.Range("E19") "  +4.24%  "

# Row 20
Set-TextValue $ws.Range("D20") "348.69"
Set-TextValue $ws.Range("E20") "  +0.27%  "

# Row 21
Set-TextValue $ws.Range("D21") "7.02"
Set-TextValue $ws.Range("E21") "  +2.68%  "

# Row 22
Set-TextValue $ws.Range("E22") "  +0.28%  "

# Row 23
Set-TextValue $ws.Range("D23") "67.38"
Set-TextValue $ws.Range("E23") "  +1.95%  "

# Row 24
Set-TextValue $ws.Range("D24") "1.69"
Set-TextValue $ws.Range("E24") "  -1.76%  "

# Row 25
Set-TextValue $ws.Range("D25") "9.25"
Set-TextValue $ws.Range("E25") "  +0.92%  "

# Row 26
Set-TextValue $ws.Range("E26") "  +0.26%  "

# Row 27
Set-TextValue $ws.Range("D27") "8.41"
Set-TextValue $ws.Range("E27") "  +4.40%  "

# Row 28
Set-TextValue $ws.Range("D28") "548.19"
Set-TextValue $ws.Range("E28") "  +0.95%  "

# Row 29
Set-TextValue $ws.Range("E29") "  -0.31%  "

# Row 30
Set-TextValue $ws.Range("E30") "  -0.06%  "

# Row 31
Set-TextValue $ws.Range("E31") "  +1.80%  "

# Row 32
Set-TextValue $ws.Range("D32") $d32val
Set-TextValue $ws.Range("E32") "  +5.71%  "

# Row 33
Set-TextValue $ws.Range("E33") "  +2.97%  "

# Row 34
Set-TextValue $ws.Range("D34") "5.43"
Set-TextValue $ws.Range("E34") "  +4.11%  "

# Row 35
Set-TextValue $ws.Range("D35") "6.16"
Set-TextValue $ws.Range("E35") "  +2.17%  "

# Row 36
Set-TextValue $ws.Range("B36") "Monero"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D36") "164.60"
Set-TextValue $ws.Range("E36") "  -2.13%  "

# Row 37
Set-TextValue $ws.Range("B37") "PolygonEcosystemToken"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue $ws.Range("D37") "0.418"
Set-TextValue $ws.Range("E37") "  +3.24%  "

# Row 38
Set-TextValue $ws.Range("E38") "  +2.27%  "

# Row 39
Set-TextValue $ws.Range("B39") "EthereumClassic"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D39") "19.81"
Set-TextValue $ws.Range("E39") "  +2.60%  "

# Row 40
Set-TextValue $ws.Range("B40") "FirstDigitalUSD"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D40") "0.999"
Set-TextValue $ws.Range("E40") "  +0.10%  "

# Row 41
Set-TextValue $ws.Range("E41") "  -0.07%  "

# Row 42
Set-TextValue $ws.Range("D42") "167.92"
Set-TextValue $ws.Range("E42") "  -0.50%  "

# Row 43
Set-TextValue $ws.Range("E43") "  +4.94%  "

# Row 44
Set-TextValue $ws.Range("D44") "23.66"
Set-TextValue $ws.Range("E44") "  +11.08%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.0586"
Set-TextValue $ws.Range("E45") "  -0.01%  "

# Row 46
Set-TextValue $ws.Range("D46") "2.16"
Set-TextValue $ws.Range("E46") "  +8.83%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.638"
Set-TextValue $ws.Range("E47") "  +1.86%  "

# Row 48
Set-TextValue $ws.Range("D48") "0.0253"
Set-TextValue $ws.Range("E48") "  +3.56%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.0967"
Set-TextValue $ws.Range("E49") "  +0.39%  "

# Row 50
Set-TextValue $ws.Range("D50") "19.27"
Set-TextValue $ws.Range("E50") "  +1.20%  "

# Row 51
Set-TextValue $ws.Range("D51") $d51val
Set-TextValue $ws.Range("E51") "  +19.44%  "
